$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Remove the "Website: www.OHare&McGovern.co.uk" paragraph and the
#    blank paragraph that follows it (paragraphs 4 and 5). Doing this
#    first keeps the paragraph indices below predictable.
# -----------------------------------------------------------------
$websitePara = $d.Paragraphs(4)
$blankPara = $d.Paragraphs(5)
$removeRange = $d.Range($websitePara.Range.Start, $blankPara.Range.End)
$removeRange.Delete()

# -----------------------------------------------------------------
# 2. Company name: "O'Hare & McGovern Ltd" -> "Creative Constructions" + " Ltd"
#    (kept as two separate runs, second with a preserved leading space)
# -----------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("O'Hare & McGovern Ltd", $true, $false, $false, $false, $false, $true, 1, $false, "Creative Constructions Ltd", 2)
$tail1 = $d.Range($rng1.Start + "Creative Constructions".Length, $rng1.End)
$tail1.Bold = 1
$tail1.Bold = 0

# -----------------------------------------------------------------
# 3. Phone number: "Phone: 08882589873" -> "Phone: " + "+44 " + "8882589873"
#    (three separate runs)
# -----------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Phone: 08882589873", $true, $false, $false, $false, $false, $true, 1, $false, "Phone: +44 8882589873", 2)
$afterPhoneLabel = $rng2.Start + "Phone: ".Length
$afterCountryCode = $rng2.Start + "Phone: +44 ".Length

$tail2a = $d.Range($afterPhoneLabel, $rng2.End)
$tail2a.Bold = 1
$tail2a.Bold = 0

$tail2b = $d.Range($afterCountryCode, $rng2.End)
$tail2b.Bold = 1
$tail2b.Bold = 0

# -----------------------------------------------------------------
# 4. "Tax Year-to-Date" -> "Tax Details"
# -----------------------------------------------------------------
$d.Content.Find.Execute("Tax Year-to-Date", $true, $false, $false, $false, $false, $true, 1, $false, "Tax Details", 2)

# -----------------------------------------------------------------
# 5. "Gross Earnings" -> "Total Income", and move the "_GoBack" bookmark
#    from the end of the "Tax Paid" paragraph to right after the
#    "Total Income" label (before the ": £12000.00" run).
# -----------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.Execute("Gross Earnings", $true, $false, $false, $false, $false, $true, 1, $false, "Total Income", 2)
$bmPoint = $d.Range($rng5.End, $rng5.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
